$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # GraphSAGE+XGBoost
$ws2 = $wb.Worksheets.Item(2)   # GraphSAGE Supervised
$ws5 = $wb.Worksheets.Item(5)   # GRAND

# ---------------------------------------------------------------------------
# Sheet2 "GraphSAGE Supervised": populate the results table.
# The order in which new string values are first written determines the
# order they are appended to sharedStrings.xml, so we write them in a
# specific sequence to reproduce the expected shared-string table.
# ---------------------------------------------------------------------------

# 1) Identifier column (rows 2-11) - introduces "gcn_small_0.0100"
$ws2.Range("B2:B11").Value = "gcn_small_0.0100"

# 2) Metric headers - introduces "Loss", "F1_Micro", "F1_Macro"
$ws2.Range("C1").Value = "Loss"
$ws2.Range("D1").Value = "F1_Micro"
$ws2.Range("E1").Value = "F1_Macro"

# 3) Test/Val column - introduces "Test" then "Val"
$ws2.Range("F2:F6").Value = "Test"
$ws2.Range("F7:F11").Value = "Val"

# 4) Test_or_Val header - introduces "Test_or_Val"
$ws2.Range("F1").Value = "Test_or_Val"

# 5) Remaining headers (reuse already-existing shared strings)
$ws2.Range("A1").Value = "Model"
$ws2.Range("B1").Value = "Identifier"

# 6) Model column data (reuse already-existing shared strings)
$ws2.Range("A2").Value = "gcn"
$ws2.Range("A3").Value = "graphsage_maxpool"
$ws2.Range("A4").Value = "graphsage_mean"
$ws2.Range("A5").Value = "graphsage_meanpool"
$ws2.Range("A6").Value = "graphsage_seq"
$ws2.Range("A7").Value = "gcn"
$ws2.Range("A8").Value = "graphsage_maxpool"
$ws2.Range("A9").Value = "graphsage_mean"
$ws2.Range("A10").Value = "graphsage_meanpool"
$ws2.Range("A11").Value = "graphsage_seq"

# 7) Numeric metric values
$ws2.Range("C2").Value = 0.30408000000000002
$ws2.Range("D2").Value = 0.85487999999999997
$ws2.Range("E2").Value = 0.48658000000000001

$ws2.Range("C3").Value = 0.23588999999999999
$ws2.Range("D3").Value = 0.88771
$ws2.Range("E3").Value = 0.70133999999999996

$ws2.Range("C4").Value = 0.24249999999999999
$ws2.Range("D4").Value = 0.89446000000000003
$ws2.Range("E4").Value = 0.75022999999999995

$ws2.Range("C5").Value = 0.21962000000000001
$ws2.Range("D5").Value = 0.90105999999999997
$ws2.Range("E5").Value = 0.75531000000000004

$ws2.Range("C6").Value = 0.23255000000000001
$ws2.Range("D6").Value = 0.90237000000000001
$ws2.Range("E6").Value = 0.77325999999999995

$ws2.Range("C7").Value = 0.32694000000000001
$ws2.Range("D7").Value = 0.85348999999999997
$ws2.Range("E7").Value = 0.48638999999999999

$ws2.Range("C8").Value = 0.2477
$ws2.Range("D8").Value = 0.88097000000000003
$ws2.Range("E8").Value = 0.69376000000000004

$ws2.Range("C9").Value = 0.24138999999999999
$ws2.Range("D9").Value = 0.90456999999999999
$ws2.Range("E9").Value = 0.77881

$ws2.Range("C10").Value = 0.21364
$ws2.Range("D10").Value = 0.91056000000000004
$ws2.Range("E10").Value = 0.78483000000000003

$ws2.Range("C11").Value = 0.22484000000000001
$ws2.Range("D11").Value = 0.90859999999999996
$ws2.Range("E11").Value = 0.78715000000000002

# 8) Header row formatting: horizontal-center-only style on A1:E1
$ws2.Range("A1:E1").HorizontalAlignment = -4108

# 9) Data cells B2:F11 reuse the workbook's existing "center both" style,
#    by copying the format from an already-styled cell on sheet1.
$ws1.Range("F2").Copy() | Out-Null
$ws2.Range("B2:F11").PasteSpecial(-4122) | Out-Null

# 10) Column widths (approximate best-fit sizing for the new content)
$ws2.Columns.Item(1).ColumnWidth = 19.5
$ws2.Columns.Item(2).ColumnWidth = 15.666666666666666
$ws2.Columns.Item(3).ColumnWidth = 7.833333333333333
$ws2.Columns.Item(4).ColumnWidth = 8.166666666666666
$ws2.Columns.Item(5).ColumnWidth = 8.666666666666666
$ws2.Columns.Item(6).ColumnWidth = 10.666666666666666

# ---------------------------------------------------------------------------
# Sheet1 "GraphSAGE+XGBoost": add the new Test_or_Val column.
# ---------------------------------------------------------------------------
$ws1.Range("G1").Value = "Test_or_Val"
$ws1.Range("G2:G7").Value = "Val"

$ws1.Range("F2").Copy() | Out-Null
$ws1.Range("G2:G7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Selections
# ---------------------------------------------------------------------------
$ws1.Range("G1").Select() | Out-Null
$ws2.Range("F13").Select() | Out-Null
$ws5.Range("E13").Select() | Out-Null

# Activate the GRAND sheet last so it becomes the active tab.
$ws5.Activate() | Out-Null
